$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 181, shifting existing rows 181-225 down to 182-226
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new record (matches row 182's
# structure/values for the constant columns, with its own date/price data)
$ws.Cells.Item(181, 1).Value = 10
$ws.Cells.Item(181, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(181, 3).Value = "La Araucanía"
$ws.Cells.Item(181, 4).Value = 44782
$ws.Cells.Item(181, 5).Value = 9
$ws.Cells.Item(181, 6).Value = 100112005
$ws.Cells.Item(181, 7).Value = "Puerro"
$ws.Cells.Item(181, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 20
$ws.Cells.Item(181, 11).Value = 16000
$ws.Cells.Item(181, 12).Value = 16000
$ws.Cells.Item(181, 13).Value = 16000
$ws.Cells.Item(181, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(181, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(181, 16).Value = 1333
$ws.Cells.Item(181, 17).Value = 12
$ws.Cells.Item(181, 18).Value = "Hortaliza"
